$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-01-09 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-01-10 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("50×49=2450", $true, $false, $false, $false, $false, $true, 1, $false, "29×24=696", 2) | Out-Null
$d.Content.Find.Execute("96×68=6528", $true, $false, $false, $false, $false, $true, 1, $false, "98×94=9212", 2) | Out-Null
$d.Content.Find.Execute("30×98=2940", $true, $false, $false, $false, $false, $true, 1, $false, "49×56=2744", 2) | Out-Null
$d.Content.Find.Execute("61×85=5185", $true, $false, $false, $false, $false, $true, 1, $false, "43×55=2365", 2) | Out-Null
$d.Content.Find.Execute("50×12=600", $true, $false, $false, $false, $false, $true, 1, $false, "80×25=2000", 2) | Out-Null
$d.Content.Find.Execute("36×57=2052", $true, $false, $false, $false, $false, $true, 1, $false, "57×58=3306", 2) | Out-Null
$d.Content.Find.Execute("66×90=5940", $true, $false, $false, $false, $false, $true, 1, $false, "29×50=1450", 2) | Out-Null
$d.Content.Find.Execute("74×96=7104", $true, $false, $false, $false, $false, $true, 1, $false, "86×32=2752", 2) | Out-Null
$d.Content.Find.Execute("16×70=1120", $true, $false, $false, $false, $false, $true, 1, $false, "54×21=1134", 2) | Out-Null
$d.Content.Find.Execute("39×36=1404", $true, $false, $false, $false, $false, $true, 1, $false, "86×45=3870", 2) | Out-Null
$d.Content.Find.Execute("88×54=4752", $true, $false, $false, $false, $false, $true, 1, $false, "13×17=221", 2) | Out-Null
$d.Content.Find.Execute("50×43=2150", $true, $false, $false, $false, $false, $true, 1, $false, "56×28=1568", 2) | Out-Null
$d.Content.Find.Execute("37×38=1406", $true, $false, $false, $false, $false, $true, 1, $false, "29×78=2262", 2) | Out-Null
$d.Content.Find.Execute("76×61=4636", $true, $false, $false, $false, $false, $true, 1, $false, "18×84=1512", 2) | Out-Null
$d.Content.Find.Execute("31×96=2976", $true, $false, $false, $false, $false, $true, 1, $false, "70×16=1120", 2) | Out-Null
$d.Content.Find.Execute("73×12=876", $true, $false, $false, $false, $false, $true, 1, $false, "45×99=4455", 2) | Out-Null
$d.Content.Find.Execute("95×89=8455", $true, $false, $false, $false, $false, $true, 1, $false, "18×45=810", 2) | Out-Null
$d.Content.Find.Execute("28×90=2520", $true, $false, $false, $false, $false, $true, 1, $false, "43×78=3354", 2) | Out-Null
$d.Content.Find.Execute("52×74=3848", $true, $false, $false, $false, $false, $true, 1, $false, "39×97=3783", 2) | Out-Null
$d.Content.Find.Execute("97×65=6305", $true, $false, $false, $false, $false, $true, 1, $false, "42×56=2352", 2) | Out-Null
$d.Content.Find.Execute("62×31=1922", $true, $false, $false, $false, $false, $true, 1, $false, "15×13=195", 2) | Out-Null
$d.Content.Find.Execute("64×87=5568", $true, $false, $false, $false, $false, $true, 1, $false, "61×46=2806", 2) | Out-Null
$d.Content.Find.Execute("80×42=3360", $true, $false, $false, $false, $false, $true, 1, $false, "21×96=2016", 2) | Out-Null
$d.Content.Find.Execute("34×16=544", $true, $false, $false, $false, $false, $true, 1, $false, "82×19=1558", 2) | Out-Null
$d.Content.Find.Execute("96×69=6624", $true, $false, $false, $false, $false, $true, 1, $false, "91×53=4823", 2) | Out-Null
